# Adjusted LooukupValue Posting due to ATTD.TestScriptor example
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

$ws.Range("F143").Value = "Posted sales invoice and shipment inherit lookup value from sales order"
$ws.Range("H144").Value = "Sales order with lookup value"
$ws.Range("H145").Value = "Post sales order (invoice & ship)"
$ws.Range("F148").Value = "Posting throws error on sales order with empty lookup value"
$ws.Range("H149").Value = "Sales order without lookup value"
$ws.Range("H150").Value = "Post sales order (invoice & ship)"
$ws.Range("F152").Value = "Posted warehouse shipment line inherits lookup value from sales order"
$ws.Range("H155").Value = "Warehouse shipment line from sales order with lookup value"
$ws.Range("H156").Value = "Post Warehouse shipment"
$ws.Range("H157").Value = "Posted warehouse shipment line has lookup value from sales order"
$ws.Range("F158").Value = "Posting throws error on warehouse shipment line with empty lookup value"
$ws.Range("H161").Value = "Warehouse shipment line from sales order without lookup value"
$ws.Range("H162").Value = "Post Warehouse shipment"
